$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 contains a duplicated "2016-17" LITO entry (mistakenly repeating
# row 5's data). Delete the entire row, shifting subsequent rows up.
$ws.Rows.Item(20).Delete()

# Update the active selection to match the fixed sheet.
$ws.Range("E20").Select()
